$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4562899804834331
$ws.Range("C2").Value = 0.1465860509869259
$ws.Range("D2").Value = 0.01678011955571179
$ws.Range("F2").Value = 0.3537062187669235
$ws.Range("G2").Value = 0.2108695015573332
$ws.Range("H2").Value = 0.3678114218550235
$ws.Range("M2").Value = 0.8957762817057073
$ws.Range("N2").Value = 0.8708063197184401
$ws.Range("O2").Value = 1.067401165052189
$ws.Range("B3").Value = 0.3981133288916965
$ws.Range("C3").Value = 0.1379728245648408
$ws.Range("D3").Value = 0.01463514096519702
$ws.Range("F3").Value = 0.3475673889692246
$ws.Range("G3").Value = 0.2059664150317246
$ws.Range("H3").Value = 0.3688040511719777
$ws.Range("M3").Value = 0.797030592983532
$ws.Range("N3").Value = 0.881288621807002
$ws.Range("O3").Value = 1.058840944304848
$ws.Range("B4").Value = 0.3622529221879347
$ws.Range("C4").Value = 0.1326524364614272
$ws.Range("D4").Value = 0.0133131152205479
$ws.Range("F4").Value = 0.3440713489052953
$ws.Range("G4").Value = 0.2031548690455836
$ws.Range("H4").Value = 0.3696172956414046
$ws.Range("M4").Value = 0.7369136299699051
$ws.Range("N4").Value = 0.8881390974238705
$ws.Range("O4").Value = 1.054424653183546
$ws.Range("B5").Value = 0.3476056276084876
$ws.Range("C5").Value = 0.1304765880226881
$ws.Range("D5").Value = 0.01277316082469326
$ws.Range("F5").Value = 0.3427153122120927
$ws.Range("G5").Value = 0.2020590061463068
$ws.Range("H5").Value = 0.3699999372135068
$ws.Range("M5").Value = 0.7125387286685054
$ws.Range("N5").Value = 0.8910349203842038
$ws.Range("O5").Value = 1.05283576999156
$ws.Range("B6").Value = 0.3451714453411796
$ws.Range("C6").Value = 0.1301148297093135
$ws.Range("D6").Value = 0.01268342938463007
$ws.Range("F6").Value = 0.3424942860682094
$ws.Range("G6").Value = 0.20188004576233
$ws.Range("H6").Value = 0.3700665691796559
$ws.Range("M6").Value = 0.7084985706507467
$ws.Range("N6").Value = 0.8915220635787833
$ws.Range("O6").Value = 1.052584657815359
$ws.Range("B7").Value = 0.3620555192144082
$ws.Range("C7").Value = 0.1326231232367405
$ws.Range("D7").Value = 0.01330583808697838
$ws.Range("F7").Value = 0.3440527831511488
$ws.Range("G7").Value = 0.2031398881769846
$ws.Range("H7").Value = 0.3696222486216385
$ws.Range("M7").Value = 0.7365844100921493
$ws.Range("N7").Value = 0.8881777295917246
$ws.Range("O7").Value = 1.054402371948896
$ws.Range("B8").Value = 0.4362603953183566
$ws.Range("C8").Value = 0.1436229690972652
$ws.Range("D8").Value = 0.0160415916299641
$ws.Range("F8").Value = 0.3515327929072072
$ws.Range("G8").Value = 0.2091375184998228
$ws.Range("H8").Value = 0.3681113875196971
$ws.Range("M8").Value = 0.8616186723733534
$ws.Range("N8").Value = 0.8743346405874632
$ws.Range("O8").Value = 1.064275092839594
$ws.Range("B9").Value = 0.5806194210595379
$ws.Range("C9").Value = 0.1649306866331415
$ws.Range("D9").Value = 0.02136523918623823
$ws.Range("F9").Value = 0.3683738152797744
$ws.Range("G9").Value = 0.222486079534832
$ws.Range("H9").Value = 0.3667658089030397
$ws.Range("M9").Value = 1.111161493213132
$ws.Range("N9").Value = 0.8504748089005574
$ws.Range("O9").Value = 1.090317445466923
$ws.Range("B10").Value = 0.6859201553729122
$ws.Range("C10").Value = 0.1804125951514095
$ws.Range("D10").Value = 0.02524979145800899
$ws.Range("F10").Value = 0.3820800932664241
$ws.Range("G10").Value = 0.2332740052372628
$ws.Range("H10").Value = 0.3667642582358326
$ws.Range("M10").Value = 1.297561531885037
$ws.Range("N10").Value = 0.8349473252965254
$ws.Range("O10").Value = 1.113554893052395
$ws.Range("B11").Value = 0.7336481379890643
$ws.Range("C11").Value = 0.1874155911477828
$ws.Range("D11").Value = 0.02701084061423131
$ws.Range("F11").Value = 0.3886068839197065
$ws.Range("G11").Value = 0.2383975662611704
$ws.Range("H11").Value = 0.3669781706024509
$ws.Range("M11").Value = 1.383116039352871
$ws.Range("N11").Value = 0.8283181222077971
$ws.Range("O11").Value = 1.125024214724903
$ws.Range("B12").Value = 0.7516953338143821
$ws.Range("C12").Value = 0.1900614723747935
$ws.Range("D12").Value = 0.02767679747920937
$ws.Range("F12").Value = 0.3911204771691033
$ws.Range("G12").Value = 0.2403690011930877
$ws.Range("H12").Value = 0.3670900501605843
$ws.Range("M12").Value = 1.415629945896626
$ws.Range("N12").Value = 0.8258702848526411
$ws.Range("O12").Value = 1.129497019370433
$ws.Range("B13").Value = 0.747809741031233
$ws.Range("C13").Value = 0.1894919058455002
$ws.Range("D13").Value = 0.02753341311166935
$ws.Range("F13").Value = 0.3905772587870899
$ws.Range("G13").Value = 0.2399430245094436
$ws.Range("H13").Value = 0.3670645814740823
$ws.Range("M13").Value = 1.408622210535825
$ws.Range("N13").Value = 0.826394690073613
$ws.Range("O13").Value = 1.128527949237252
$ws.Range("B14").Value = 0.7351334276168586
$ws.Range("C14").Value = 0.1876333912559858
$ws.Range("D14").Value = 0.0270656478597715
$ws.Range("F14").Value = 0.3888128357180136
$ws.Range("G14").Value = 0.2385591299929928
$ws.Range("H14").Value = 0.3669867561285542
$ws.Range("M14").Value = 1.385788607283686
$ws.Range("N14").Value = 0.8281154843414669
$ws.Range("O14").Value = 1.12538959484732
$ws.Range("B15").Value = 0.727365347151931
$ws.Range("C15").Value = 0.186494207606728
$ws.Range("D15").Value = 0.02677900777486997
$ws.Range("F15").Value = 0.3877375529953255
$ws.Range("G15").Value = 0.2377155301605427
$ws.Range("H15").Value = 0.3669431071829621
$ws.Range("M15").Value = 1.371817717419873
$ws.Range("N15").Value = 0.8291776609478632
$ws.Range("O15").Value = 1.12348415514478
$ws.Range("B16").Value = 0.6827974052305876
$ws.Range("C16").Value = 0.1799541082153553
$ws.Range("D16").Value = 0.02513457720902323
$ws.Range("F16").Value = 0.3816594301615339
$ws.Range("G16").Value = 0.2329435320819044
$ws.Range("H16").Value = 0.3667545983951612
$ws.Range("M16").Value = 1.291986231332586
$ws.Range("N16").Value = 0.8353893025650407
$ws.Range("O16").Value = 1.11282346171825
$ws.Range("B17").Value = 0.6554109242825916
$ws.Range("C17").Value = 0.175931577042121
$ws.Range("D17").Value = 0.02412419083734818
$ws.Range("F17").Value = 0.3780054821158387
$ws.Range("G17").Value = 0.2300715383415763
$ws.Range("H17").Value = 0.3666939328959984
$ws.Range("M17").Value = 1.243211798953737
$ws.Range("N17").Value = 0.8393112133709906
$ws.Range("O17").Value = 1.106513895326373
$ws.Range("B18").Value = 0.6396426622316937
$ws.Range("C18").Value = 0.1736141937588798
$ws.Range("D18").Value = 0.02354247619299343
$ws.Range("F18").Value = 0.3759312822167544
$ws.Range("G18").Value = 0.2284399711133887
$ws.Range("H18").Value = 0.3666792397170013
$ws.Range("M18").Value = 1.215229162080149
$ws.Range("N18").Value = 0.8416078707879748
$ws.Range("O18").Value = 1.102969344761902
$ws.Range("B19").Value = 0.6343010444263086
$ws.Range("C19").Value = 0.17282893548267
$ws.Range("D19").Value = 0.02334542153770514
$ws.Range("F19").Value = 0.3752337063899773
$ws.Range("G19").Value = 0.2278910362272057
$ws.Range("H19").Value = 0.3666777338099223
$ws.Range("M19").Value = 1.205766729017441
$ws.Range("N19").Value = 0.8423924991671399
$ws.Range("O19").Value = 1.101783728746398
$ws.Range("B20").Value = 0.658327956938308
$ws.Range("C20").Value = 0.1763601703812014
$ws.Range("D20").Value = 0.02423180720340667
$ws.Range("F20").Value = 0.378391609060543
$ws.Range("G20").Value = 0.2303751618311622
$ws.Range("H20").Value = 0.3666983001184718
$ws.Range("M20").Value = 1.24839650937912
$ws.Range("N20").Value = 0.8388894882607119
$ws.Range("O20").Value = 1.107176806298696
$ws.Range("B21").Value = 0.7388574942246464
$ws.Range("C21").Value = 0.1881794471310911
$ws.Range("D21").Value = 0.02720306701489505
$ws.Range("F21").Value = 0.3893299479120742
$ws.Range("G21").Value = 0.2389647637010199
$ws.Range("H21").Value = 0.3670087772532895
$ws.Range("M21").Value = 1.39249217817077
$ws.Range("N21").Value = 0.8276083488095693
$ws.Range("O21").Value = 1.126307884377013
$ws.Range("B22").Value = 0.7913339541170785
$ws.Range("C22").Value = 0.1958689183335593
$ws.Range("D22").Value = 0.0291396078243622
$ws.Range("F22").Value = 0.3967238620839524
$ws.Range("G22").Value = 0.2447608187170403
$ws.Range("H22").Value = 0.3673916650842841
$ws.Range("M22").Value = 1.487348101450095
$ws.Range("N22").Value = 0.8205997816700474
$ws.Range("O22").Value = 1.139566877046434
$ws.Range("B23").Value = 0.7633408414161522
$ws.Range("C23").Value = 0.1917682071162403
$ws.Range("D23").Value = 0.02810654350477648
$ws.Range("F23").Value = 0.3927551374498748
$ws.Range("G23").Value = 0.2416506189734235
$ws.Range("H23").Value = 0.3671708384857624
$ws.Range("M23").Value = 1.4366570470808
$ws.Range("N23").Value = 0.8243070360683618
$ws.Range("O23").Value = 1.132421015553007
$ws.Range("B24").Value = 0.6570092392864808
$ws.Range("C24").Value = 0.1761664181888705
$ws.Range("D24").Value = 0.0241831564266235
$ws.Range("F24").Value = 0.3782169585229767
$ws.Range("G24").Value = 0.2302378326774743
$ws.Range("H24").Value = 0.3666962628288104
$ws.Range("M24").Value = 1.246052320310824
$ws.Range("N24").Value = 0.8390800196811838
$ws.Range("O24").Value = 1.10687684605486
$ws.Range("B25").Value = 0.5416959624907065
$ws.Range("C25").Value = 0.1591958358306584
$ws.Range("D25").Value = 0.01992962688948552
$ws.Range("F25").Value = 0.3635843719167298
$ws.Range("G25").Value = 0.2187036582495665
$ws.Range("H25").Value = 0.3669565714952654
$ws.Range("M25").Value = 1.043146084565535
$ws.Range("N25").Value = 0.8565778325861473
$ws.Range("O25").Value = 1.082553779536084
